$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.084595918655396
$ws.Range("B1").Value = 2.285298585891724
$ws.Range("C1").Value = 2.500075340270996
$ws.Range("D1").Value = 3.654590606689453
$ws.Range("E1").Value = 1.661850452423096
